$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "weird typo" Propionate value in L10 -> mark as NA like its neighbours
$ws.Range("L10").Value = "NA"

# Fill in the new Substrate / Antibiotics summary rows (12 and 13)
$ws.Range("A12").Value = "Substrate"
$ws.Range("A13").Value = "Antibiotics"

$ws.Range("B12:O12").Value = "Inulin"
$ws.Range("P12:V12").Value = "Saline"

$ws.Range("B13:H13").Value = "No"
$ws.Range("I13:O13").Value = "Yes"
$ws.Range("P13:S13").Value = "No"
$ws.Range("T13:V13").Value = "Yes"

# J12 previously held a leftover empty-cell style; align it with the rest of
# row 12's default column formatting by copying format from a sibling cell.
$ws.Range("I12").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author's last edit landed
$ws.Range("V13").Select()

Write-Output "done"
